$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 12827051
$ws.Range("I76").Value = 7714.45
$ws.Range("J76").Value = 55558172
$ws.Range("K76").Value = 7714.45
$ws.Range("L76").Value = 55558172
$ws.Range("M76").Value = -7399.45
$ws.Range("N76").Value = -55558802

$ws.Range("H79").Value = 12827051
$ws.Range("I79").Value = 7714.45
$ws.Range("J79").Value = 55558172
$ws.Range("K79").Value = 7714.45
$ws.Range("L79").Value = 55558172
$ws.Range("M79").Value = -6622.45
$ws.Range("N79").Value = -55560356

$ws.Range("H107").Value = 700
$ws.Range("I107").Value = 350
$ws.Range("J107").Value = 933.3333
$ws.Range("K107").Value = 350
$ws.Range("L107").Value = 933.3333
$ws.Range("M107").Value = 1570
$ws.Range("N107").Value = -4773.3333

$ws.Range("H112").Value = 15550.454
$ws.Range("J112").Value = 16732.844
$ws.Range("L112").Value = 50198.53200000001
$ws.Range("N112").Value = -52414.53200000001

$ws.Range("H125").Value = 842.4286
$ws.Range("I125").Value = 799
$ws.Range("J125").Value = 875
$ws.Range("K125").Value = 7191
$ws.Range("L125").Value = 7875
$ws.Range("M125").Value = -4731
$ws.Range("N125").Value = -12795

$ws.Range("H129").Value = 1039.2285
$ws.Range("J129").Value = 925.6875
$ws.Range("L129").Value = 2777.0625
$ws.Range("N129").Value = -12777.0625

$ws.Range("H132").Value = 7580676
$ws.Range("I132").Value = 10004451
$ws.Range("J132").Value = 6381.0625
$ws.Range("K132").Value = 30013353
$ws.Range("L132").Value = 19143.1875
$ws.Range("M132").Value = -30010823
$ws.Range("N132").Value = -24203.1875

$ws.Range("H135").Value = 832.6957
$ws.Range("I135").Value = 712
$ws.Range("J135").Value = 2100
$ws.Range("K135").Value = 6408
$ws.Range("L135").Value = 18900
$ws.Range("M135").Value = -3873
$ws.Range("N135").Value = -23970

$ws.Range("H137").Value = 1399.5641
$ws.Range("I137").Value = 1220.8214
$ws.Range("J137").Value = 1854.5454
$ws.Range("K137").Value = 3662.4642
$ws.Range("L137").Value = 5563.6362
$ws.Range("M137").Value = -1112.4642
$ws.Range("N137").Value = -10663.6362

$ws.Range("H138").Value = 1947.1803
$ws.Range("I138").Value = 724.13336
$ws.Range("J138").Value = 3130.7742
$ws.Range("K138").Value = 2172.40008
$ws.Range("L138").Value = 9392.3226
$ws.Range("M138").Value = 2967.59992
$ws.Range("N138").Value = -19672.3226

$ws.Range("H141").Value = 3025.5278
$ws.Range("I141").Value = 1620.3
$ws.Range("J141").Value = 10051.667
$ws.Range("K141").Value = 4860.9
$ws.Range("L141").Value = 30155.001
$ws.Range("M141").Value = 319.1000000000004
$ws.Range("N141").Value = -40515.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 807.1064
$ws.Range("I74").Value = 788.9048
$ws.Range("J74").Value = 960
$ws.Range("K74").Value = 788.9048
$ws.Range("L74").Value = 960
$ws.Range("M74").Value = 85.09519999999998
$ws.Range("N74").Value = -2708

$ws.Range("H77").Value = 807.1064
$ws.Range("I77").Value = 788.9048
$ws.Range("J77").Value = 960
$ws.Range("K77").Value = 3944.524
$ws.Range("L77").Value = 4800
$ws.Range("M77").Value = 423.4759999999997
$ws.Range("N77").Value = -13536

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2367143
$ws.Range("I134").Value = 794.84375
$ws.Range("J134").Value = 7415352.5
$ws.Range("K134").Value = 2384.53125
$ws.Range("L134").Value = 22246057.5
$ws.Range("M134").Value = 150.46875
$ws.Range("N134").Value = -22251127.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 40000
$ws.Range("J20").Value = 40000
$ws.Range("L20").Value = 40000
$ws.Range("N20").Value = -40472

$ws.Range("H30").Value = 40000
$ws.Range("J30").Value = 40000
$ws.Range("L30").Value = 40000
$ws.Range("N30").Value = -40182

$ws.Range("H31").Value = 1163.0615
$ws.Range("I31").Value = 862.1
$ws.Range("J31").Value = 1644.6
$ws.Range("K31").Value = 862.1
$ws.Range("L31").Value = 1644.6
$ws.Range("M31").Value = -567.1
$ws.Range("N31").Value = -2234.6

$ws.Range("H34").Value = 1163.0615
$ws.Range("I34").Value = 862.1
$ws.Range("J34").Value = 1644.6
$ws.Range("K34").Value = 862.1
$ws.Range("L34").Value = 1644.6
$ws.Range("M34").Value = -660.1
$ws.Range("N34").Value = -2048.6

$ws.Range("H58").Value = 18868854
$ws.Range("I58").Value = 25642038
$ws.Range("J58").Value = 701
$ws.Range("K58").Value = 25642038
$ws.Range("L58").Value = 701
$ws.Range("M58").Value = -25641835
$ws.Range("N58").Value = -1107

$ws.Range("H128").Value = 40000
$ws.Range("J128").Value = 40000
$ws.Range("L128").Value = 40000
$ws.Range("N128").Value = -49960

$ws.Range("H132").Value = 10102441
$ws.Range("I132").Value = 1284.7693
$ws.Range("J132").Value = 47621020
$ws.Range("K132").Value = 3854.3079
$ws.Range("L132").Value = 142863060
$ws.Range("M132").Value = -1324.3079
$ws.Range("N132").Value = -142868120

$ws.Range("H134").Value = 804.28
$ws.Range("I134").Value = 849.04877
$ws.Range("J134").Value = 600.3333
$ws.Range("K134").Value = 2547.14631
$ws.Range("L134").Value = 1800.9999
$ws.Range("M134").Value = -12.14631000000008
$ws.Range("N134").Value = -6870.9999

$ws.Range("H136").Value = 18868854
$ws.Range("I136").Value = 25642038
$ws.Range("J136").Value = 701
$ws.Range("K136").Value = 76926114
$ws.Range("L136").Value = 2103
$ws.Range("M136").Value = -76923564
$ws.Range("N136").Value = -7203

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H126").Value = 20836034
$ws.Range("I126").Value = 166667260
$ws.Range("K126").Value = 500001780
$ws.Range("M126").Value = -499996840

$ws.Range("H131").Value = 3493.3645
$ws.Range("I131").Value = 1466.9
$ws.Range("J131").Value = 3729
$ws.Range("K131").Value = 4400.700000000001
$ws.Range("L131").Value = 11187
$ws.Range("M131").Value = 639.2999999999993
$ws.Range("N131").Value = -21267

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 1206.5
$ws.Range("I9").Value = 1206.5
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 1206.5
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -1036.5
$ws.Range("N9").ClearContents()

$ws.Range("H100").Value = 37975
$ws.Range("J100").Value = 37975
$ws.Range("L100").Value = 37975
$ws.Range("N100").Value = -40139

$ws.Range("H132").Value = 5201.523
$ws.Range("I132").Value = 3166.675
$ws.Range("J132").Value = 25550
$ws.Range("K132").Value = 9500.025000000001
$ws.Range("L132").Value = 76650
$ws.Range("M132").Value = -6970.025000000001
$ws.Range("N132").Value = -81710

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2152
$ws.Range("I7").Value = 2003.3334
$ws.Range("J7").Value = 2375
$ws.Range("K7").Value = 2003.3334
$ws.Range("L7").Value = 2375
$ws.Range("M7").Value = -1891.3334
$ws.Range("N7").Value = -2599

$ws.Range("H22").Value = 1957.6666
$ws.Range("I22").Value = 401
$ws.Range("J22").Value = 2269
$ws.Range("K22").Value = 401
$ws.Range("L22").Value = 2269
$ws.Range("M22").Value = -106
$ws.Range("N22").Value = -2859

$ws.Range("H27").Value = 1957.6666
$ws.Range("I27").Value = 401
$ws.Range("J27").Value = 2269
$ws.Range("K27").Value = 401
$ws.Range("L27").Value = 2269
$ws.Range("M27").Value = -294
$ws.Range("N27").Value = -2483

$ws.Range("H126").Value = 2152
$ws.Range("I126").Value = 2003.3334
$ws.Range("J126").Value = 2375
$ws.Range("K126").Value = 6010.0002
$ws.Range("L126").Value = 7125
$ws.Range("M126").Value = -3540.0002
$ws.Range("N126").Value = -12065

$ws.Range("H132").Value = 31258654
$ws.Range("I132").Value = 38463690
$ws.Range("J132").Value = 36834.168
$ws.Range("K132").Value = 115391070
$ws.Range("L132").Value = 110502.504
$ws.Range("M132").Value = -115388540
$ws.Range("N132").Value = -115562.504

$ws.Range("H136").Value = 35716416
$ws.Range("I136").Value = 5293253
$ws.Range("J136").Value = 200001500
$ws.Range("K136").Value = 15879759
$ws.Range("L136").Value = 600004500
$ws.Range("M136").Value = -15877209
$ws.Range("N136").Value = -600009600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 15986.629
$ws.Range("I122").Value = 23969.727
$ws.Range("J122").Value = 2476.7693
$ws.Range("K122").Value = 71909.181
$ws.Range("L122").Value = 7430.3079
$ws.Range("M122").Value = -69459.181
$ws.Range("N122").Value = -12330.3079

$ws.Range("H126").Value = 1275.2307
$ws.Range("I126").Value = 827.8
$ws.Range("J126").Value = 2766.6667
$ws.Range("K126").Value = 2483.4
$ws.Range("L126").Value = 8300.000100000001
$ws.Range("M126").Value = -13.39999999999964
$ws.Range("N126").Value = -13240.0001

$ws.Range("H132").Value = 20916.176
$ws.Range("I132").Value = 21739.02
$ws.Range("J132").Value = 15876.25
$ws.Range("K132").Value = 65217.06
$ws.Range("L132").Value = 47628.75
$ws.Range("M132").Value = -62687.06
$ws.Range("N132").Value = -52688.75

$ws.Range("H136").Value = 11633440
$ws.Range("I136").Value = 16135584
$ws.Range("J136").Value = 2904.1667
$ws.Range("K136").Value = 48406752
$ws.Range("L136").Value = 8712.500100000001
$ws.Range("M136").Value = -48404202
$ws.Range("N136").Value = -13812.5001

$ws.Range("H137").Value = 48140
$ws.Range("J137").Value = 48140
$ws.Range("L137").Value = 48140
$ws.Range("N137").Value = -58340

Write-Output "applied updates"
